# parameterPID.xlsx - add "improved for closed loop" note + clone the P/I/D
# parameter matrix into a second copy below it (tweaking the I row's D
# value), and recolor the workbook's background (window) theme color.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New label row introducing the cloned matrix.
$ws.Range("B11").Value = "imporved for close loop"

# Clone the existing P/I/D parameter matrix (B1:E5) down to B12:E16 -
# Copy/PasteSpecial carries over both values and number formats/styles.
$ws.Range("B1:E5").Copy()
$ws.Range("B12").PasteSpecial()

# Make sure the scientific-notation formatting on the D-row survives
# the clone with the exact style used by the original matrix.
$ws.Range("C16:E16").NumberFormat = "0.00E+00"

# Tweak the cloned I row's speed(Hz) value.
$ws.Range("D15").Value = 0.3

# Leave the selection where Excel would after pasting over the new rows.
[void]$ws.Range("B12").Select()

# Recolor the theme's window/background color (Page Layout > Colors >
# Customize Colors > Background 1).
$tcs = $wb.Theme.ThemeColorScheme
$tcs.Colors(2).RGB = 13430215
